$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.119, -0.401, 0.264, -74, 39, -34),
    @(0.001, 0.124, -0.064, 10, -68, -13),
    @(0.019, -0.297, 0.887, 90, 12, 58),
    @(-0.022, -0.819, 0.211, 172, 10, 134),
    @(-0.162, -0.26, -0.289, 97, 33, 61),
    @(-0.082, -0.134, 0.6899999999999999, -116, 57, -74),
    @(-0.062, -0.043, 0.18, 99, -61, 153),
    @(-0.156, 0.037, -0.007, -79, -39, 173),
    @(0.49, -0.112, -0.262, -22, -11, 5),
    @(0.31, 0.192, 0.083, -85, 29, -116)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    for ($c = 0; $c -lt $data[$r].Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $data[$r][$c]
    }
}
